# Auto-generated edit script for cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.463.74"
$ws.Range("E2").Value = "  -0.95%  "

# Row 3
$ws.Range("D3").Value = "1.591.16"
$ws.Range("E3").Value = "  -0.59%  "

# Row 4
$ws.Range("E4").Value = "  +0.22%  "

# Row 5
$ws.Range("D5").Value = "209.81"
$ws.Range("E5").Value = "  -0.84%  "

# Row 6
$ws.Range("E6").Value = "  -1.02%  "

# Row 7
$ws.Range("E7").Value = "  +0.25%  "

# Row 8
$ws.Range("E8").Value = "  -1.03%  "

# Row 9
$ws.Range("D9").Value = "0.244"
$ws.Range("E9").Value = "  -1.33%  "

# Row 10
$ws.Range("D10").Value = "19.48"
$ws.Range("E10").Value = "  -0.45%  "

# Row 11
$ws.Range("D11").Value = "0.0844"
$ws.Range("E11").Value = "  +0.19%  "

# Row 12
$ws.Range("D12").Value = "1.815.62"
$ws.Range("E12").Value = "  -0.51%  "

# Row 13
$ws.Range("D13").Value = "1.595.20"
$ws.Range("E13").Value = "  -0.53%  "

# Row 14
$ws.Range("D14").Value = "4.02"
$ws.Range("E14").Value = "  -0.95%  "

# Row 15
$ws.Range("D15").Value = "0.516"
$ws.Range("E15").Value = "  -1.27%  "

# Row 16
$ws.Range("D16").Value = "64.09"
$ws.Range("E16").Value = "  -1.94%  "

# Row 17
$ws.Range("E17").Value = "  -2.65%  "

# Row 18
$ws.Range("E18").Value = "  +0.24%  "

# Row 19
$ws.Range("D19").Value = "206.37"
$ws.Range("E19").Value = "  -1.30%  "

# Row 20
$ws.Range("D20").Value = "7.04"
$ws.Range("E20").Value = "  -2.55%  "

# Row 21
$ws.Range("E21").Value = "  -0.92%  "

# Row 22
$ws.Range("D22").Value = "'2.20"
$ws.Range("E22").Value = "  -4.69%  "

# Row 23
$ws.Range("D23").Value = "8.86"
$ws.Range("E23").Value = "  -1.01%  "

# Row 24
$ws.Range("D24").Value = "144.57"
$ws.Range("E24").Value = "  +1.46%  "

# Row 25
$ws.Range("E25").Value = "  +0.28%  "

# Row 26
$ws.Range("D26").Value = "7.03"
$ws.Range("E26").Value = "  -1.03%  "

# Row 27
$ws.Range("E27").Value = "  -1.53%  "

# Row 28
$ws.Range("D28").Value = "15.19"
$ws.Range("E28").Value = "  -1.07%  "

# Row 29
$ws.Range("D29").Value = "0.0502"
$ws.Range("E29").Value = "  -3.83%  "

# Row 30
$ws.Range("E30").Value = "  -1.10%  "

# Row 31
$ws.Range("E31").Value = "  -0.89%  "

# Row 32
$ws.Range("E32").Value = "  -1.37%  "

# Row 33
$ws.Range("D33").Value = "1.279.10"
$ws.Range("E33").Value = "  -1.17%  "

# Row 34
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "2.49"
$ws.Range("E34").Value = "  +1.09%  "

# Row 35
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "1.24"
$ws.Range("E35").Value = "  +13.18%  "

# Row 36
$ws.Range("D36").Value = "0.595"
$ws.Range("E36").Value = "  -4.66%  "

# Row 37
$ws.Range("E37").Value = "  -2.31%  "

# Row 38
$ws.Range("D38").Value = "0.0166"
$ws.Range("E38").Value = "  -2.57%  "

# Row 39
$ws.Range("D39").Value = "0.813"
$ws.Range("E39").Value = "  -1.35%  "

# Row 40
$ws.Range("D40").Value = "5.38"
$ws.Range("E40").Value = "  -0.94%  "

# Row 41
$ws.Range("D41").Value = "'2.20"
$ws.Range("E41").Value = "  -0.09%  "

# Row 42
$ws.Range("E42").Value = "  -1.97%  "

# Row 43
$ws.Range("D43").Value = "62.13"
$ws.Range("E43").Value = "  -1.69%  "

# Row 44
$ws.Range("D44").Value = "1.726.71"
$ws.Range("E44").Value = "  -0.55%  "

# Row 45
$ws.Range("D45").Value = "88.69"
$ws.Range("E45").Value = "  -2.92%  "

# Row 46
$ws.Range("E46").Value = "  -1.26%  "

# Row 47
$ws.Range("E47").Value = "  +0.85%  "

# Row 48
$ws.Range("D48").Value = "{0}{1}{2}" -f "0.0", [char]0x2086, "0101"
$ws.Range("E48").Value = "  -4.76%  "

# Row 49
$ws.Range("E49").Value = "  +0.14%  "

# Row 50
$ws.Range("E50").Value = "  +0.23%  "

# Row 51
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "'0.400"
$ws.Range("E51").Value = "  +1.64%  "
